$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.143.17'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').Value = '1.789.18'
$ws.Range('E3').Value = '  -1.29%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.25'
$ws.Range('D5').ClearFormats()

$ws.Range('E6').Value = '  +1.98%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.69'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.64%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.17'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.87%  '

$ws.Range('E10').Value = '  +1.32%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0663'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.68%  '

$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').Value = '2.046.29'
$ws.Range('E13').Value = '  -1.29%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.45'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +12.82%  '

$ws.Range('D15').Value = '1.790.05'
$ws.Range('E15').Value = '  -1.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.633'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.71%  '

$ws.Range('D17').Value = '34.132.52'
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.24'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.54%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.55'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.68%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '254.98'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.13%  '

$ws.Range('D21').Value = '0.0₃0745'
$ws.Range('E21').Value = '  +0.42%  '

$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.52'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.44%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.16'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.59'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.34%  '

$ws.Range('E28').Value = '  +0.18%  '

$ws.Range('E29').Value = '  -0.40%  '

$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.82'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.36%  '

$ws.Range('E32').Value = '  +1.77%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.59'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.47%  '

$ws.Range('E35').Value = '  +2.53%  '

$ws.Range('D36').Value = '1.455.27'
$ws.Range('E36').Value = '  -6.01%  '

$ws.Range('E37').Value = '  -0.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.636'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.95%  '

$ws.Range('E39').Value = '  +0.89%  '

$ws.Range('E40').Value = '  +2.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '83.58'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.48%  '

$ws.Range('E42').Value = '  +0.26%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.903'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('E44').Value = '  -0.45%  '

$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.945.85'
$ws.Range('E47').Value = '  -1.00%  '

$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.84'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.24%  '

$ws.Range('E49').Value = '  +8.28%  '

$ws.Range('E50').Value = '  +0.05%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.14'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.82%  '
